$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pascal Siakam (row 7) gains "C" eligibility: SF,PF -> SF,PF,C
$ws.Range("B7").Value = "SF,PF,C"

# Ayo Dosunmu (row 15) gains "PG" eligibility: SG,SF -> PG,SG,SF
$ws.Range("B15").Value = "PG,SG,SF"
